# Update "想去人数" (want-to-go count) values in column F for the
# "展览" (rId1) and "全部类型" (rId4) sheets, reflecting a refreshed
# data scrape (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 15840
$ws1.Range("F7").Value = 10
$ws1.Range("F8").Value = 710
$ws1.Range("F9").Value = 15473
$ws1.Range("F10").Value = 59
$ws1.Range("F11").Value = 9087
$ws1.Range("F14").Value = 1017
$ws1.Range("F15").Value = 105
$ws1.Range("F16").Value = 205
$ws1.Range("F18").Value = 206
$ws1.Range("F20").Value = 63
$ws1.Range("F21").Value = 566
$ws1.Range("F23").Value = 12
$ws1.Range("F24").Value = 63
$ws1.Range("F25").Value = 1118
$ws1.Range("F26").Value = 5
$ws1.Range("F32").Value = 411
$ws1.Range("F34").Value = 48
$ws1.Range("F36").Value = 328
$ws1.Range("F37").Value = 462
$ws1.Range("F39").Value = 5584
$ws1.Range("F40").Value = 5231

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 15841
$ws4.Range("F7").Value = 10
$ws4.Range("F8").Value = 710
$ws4.Range("F9").Value = 15473
$ws4.Range("F10").Value = 59
$ws4.Range("F11").Value = 9087
$ws4.Range("F14").Value = 1017
$ws4.Range("F15").Value = 105
$ws4.Range("F16").Value = 205
$ws4.Range("F18").Value = 206
$ws4.Range("F20").Value = 63
$ws4.Range("F21").Value = 566
$ws4.Range("F23").Value = 12
$ws4.Range("F24").Value = 63
$ws4.Range("F25").Value = 1118
$ws4.Range("F26").Value = 5
$ws4.Range("F34").Value = 411
$ws4.Range("F36").Value = 48
$ws4.Range("F38").Value = 328
$ws4.Range("F39").Value = 462
$ws4.Range("F41").Value = 5584
$ws4.Range("F43").Value = 5231
